$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 431
$ws.Range("F3").Value = 2747
$ws.Range("F9").Value = 591
$ws.Range("F12").Value = 11432
$ws.Range("F13").Value = 6502
$ws.Range("F16").Value = 410
$ws.Range("F20").Value = 903
$ws.Range("F21").Value = 40
$ws.Range("F22").Value = 253
$ws.Range("F24").Value = 3626
$ws.Range("F27").Value = 493
$ws.Range("F29").Value = 309
$ws.Range("F31").Value = 289
$ws.Range("F32").Value = 4991
$ws.Range("F34").Value = 1225
$ws.Range("F35").Value = 218
$ws.Range("F36").Value = 406
$ws.Range("F37").Value = 178

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 3658
$ws.Range("F12").Value = 85

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8987
$ws.Range("F3").Value = 483
$ws.Range("F4").Value = 1791

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8987
$ws.Range("F3").Value = 483
$ws.Range("F4").Value = 1791
$ws.Range("F5").Value = 431
$ws.Range("F6").Value = 2747
$ws.Range("F16").Value = 591
$ws.Range("F19").Value = 11432
$ws.Range("F20").Value = 3658
$ws.Range("F21").Value = 6502
$ws.Range("F22").Value = 85
$ws.Range("F28").Value = 903
$ws.Range("F29").Value = 40
$ws.Range("F31").Value = 3626
$ws.Range("F34").Value = 309
$ws.Range("F39").Value = 4991
$ws.Range("F41").Value = 1225
$ws.Range("F43").Value = 218
$ws.Range("F44").Value = 178
